$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A113").Value = "Table_CaptureRatios"
$ws.Range("B113").Value = "Test Table_CaptureRatios"
$ws.Range("C113").Value = "Table_CaptureRatios_test"

$ws.Range("A114").Value = "Table_UpDownRatios"
$ws.Range("B114").Value = "Test Table_UpDownRatios"
$ws.Range("C114").Value = "Table_UpDownTatios_test"

$ws.Range("F111").Select()
